$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 98 (pushes the existing "Y1" crystal row down to row 99)
$ws.Rows("98:98").Insert()

# Populate the new row 98 with the XBee SMD header part data
$ws.Range("A98").Value = "X3"
$ws.Range("B98").Value = "XBEE-1B3"
$ws.Range("C98").Value = "XBEE-SMD"
$ws.Range("D98").Value = "XBEE-SMD"
$ws.Range("E98").Value = "SFE"
$ws.Range("F98").Value = "PRT-10030"
$ws.Range("G98").Value = "NPPN101BFLD-RC"
$ws.Range("H98").Value = 2
$ws.Range("I98").Value = 0.95
$ws.Range("J98").Formula = "=H98*I98"

# Restore the selection to match the new active cell location
[void]$ws.Range("J98").Select()

# Best-effort: scroll the view down so row 98 area is visible (topLeftCell ~ A76)
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
